$wb = $excel.ActiveWorkbook

# Sheet "展览" (rId1 / sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 964
$ws1.Range("F3").Value = 1918
$ws1.Range("F4").Value = 426

# Sheet "全部类型" (rId4 / sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 964
$ws4.Range("F5").Value = 1918
$ws4.Range("F6").Value = 426
